$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11 (shifts existing rows 11..72 down to 12..73)
$ws.Rows.Item(11).Insert()

# Fix two existing values (row numbers are now in their pre-shift positions
# since rows 2 and 10 are above the inserted row and are unaffected by the insert)
$ws.Cells.Item(2, 2).Value = 18005.48
$ws.Cells.Item(10, 2).Value = 14116.1

# Populate the newly inserted row 11 with the new data point
$ws.Cells.Item(11, 1).Value = 14
$ws.Cells.Item(11, 2).Value = 17897.95
$ws.Cells.Item(11, 3).Value = 7
$ws.Cells.Item(11, 4).Value = 2025
$ws.Cells.Item(11, 5).Value = "07/2025"
